# WageLedgerReport header: replace the U-column formula-reference plumbing
# with plain text header cells, merge the info cells, and restyle with the
# IPAPGothic font + alignment used by the new header layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Department / Employee / Sex header cells (row 4) ------------------
# Labels
$ws.Range("B4").Value = "DepartmentLabel"
$ws.Range("G4").Value = "EmployeeLabel"
$ws.Range("L4").Value = "SexLabel"

# Info cells (previously formulas referencing the U column) become plain text
$ws.Range("C4").Value = "DepartmentInfo"
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

$ws.Range("H4").Value = "EmployeeInfo"
$ws.Range("J4").ClearContents()

$ws.Range("M4").Value = "SexInfo"

# --- 2. Merge the info cells ----------------------------------------------
$ws.Range("C4:E4").Merge()
$ws.Range("H4:J4").Merge()

# --- 3. Clear the old helper cells (U1:U5) that fed the formulas ----------
$ws.Range("U1").ClearContents()
$ws.Range("U2").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("U5").ClearContents()

# --- 4. Fonts: introduce IPAPGothic in place of Calibri --------------------
# Label cells: right aligned, vertically centered
$ws.Range("B4").Font.Name = "IPAPGothic"
$ws.Range("B4").HorizontalAlignment = -4152   # xlRight
$ws.Range("B4").VerticalAlignment = -4108     # xlCenter

$ws.Range("L4").Font.Name = "IPAPGothic"
$ws.Range("L4").HorizontalAlignment = -4152   # xlRight
$ws.Range("L4").VerticalAlignment = -4108     # xlCenter

# Info cells: left aligned, vertically centered
$ws.Range("C4:E4").Font.Name = "IPAPGothic"
$ws.Range("C4:E4").HorizontalAlignment = -4131  # xlLeft
$ws.Range("C4:E4").VerticalAlignment = -4108    # xlCenter

$ws.Range("H4:J4").Font.Name = "IPAPGothic"
$ws.Range("H4:J4").HorizontalAlignment = -4131  # xlLeft
$ws.Range("H4:J4").VerticalAlignment = -4108    # xlCenter

$ws.Range("M4").Font.Name = "IPAPGothic"
$ws.Range("M4").HorizontalAlignment = -4131     # xlLeft
$ws.Range("M4").VerticalAlignment = -4108       # xlCenter

# EmployeeLabel cell keeps the new font but no explicit alignment override
$ws.Range("G4").Font.Name = "IPAPGothic"

# The gray placeholder cells (T,U,V,W,X columns) switch font too
$ws.Range("T1:T6").Font.Name = "IPAPGothic"
$ws.Range("U1:U5").Font.Name = "IPAPGothic"
$ws.Range("V1:X6").Font.Name = "IPAPGothic"

Write-Host "WageLedgerReport header applied"
